# Week 17 data logging + tiebreaking fix
# This script appends the Week 17 per-game numeric samples to the existing
# shared "history" strings on the YDS and ST sheets, and updates the
# season-to-date aggregate totals on OFF, DEF, ST, TURNS and PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append Week 17 per-play yardage samples
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Text + " 2 -2 6 6 2 13 5 1 1 13 16 6 7 -1 0 2 2 6 3 6 1 10 3 6 6 1 8 2 -2 0 37 3 4 5 2 2"
$ws.Range("B3").Value = $ws.Range("B3").Text + " 3 3 12 12 12 1 6 17 27 12 43 6 8 13 45 29"
$ws.Range("C2").Value = $ws.Range("C2").Text + " -4 6 3 -1 2 6 4 2 5 3 6 5 -1 2 6 4 2 4 3 2 8 2 8 -1 5 3 2"
$ws.Range("C3").Value = $ws.Range("C3").Text + " 6 3 9 -2 2 8 6 15 24 0 8 6 9 4 8 3 20 9 5 6 14"

# ---------------------------------------------------------------------
# OFF sheet: season-to-date offensive totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 452
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 127
$ws.Range("G2").Value = 128
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 58
$ws.Range("L2").Value = 544
$ws.Range("M2").Value = 349
$ws.Range("O2").Value = 36
$ws.Range("Q2").Value = 1017
$ws.Range("C3").Value = 317
$ws.Range("E3").Value = 85
$ws.Range("F3").Value = 191
$ws.Range("H3").Value = 56
$ws.Range("I3").Value = 112
$ws.Range("J3").Value = 83
$ws.Range("N3").Value = 39

# ---------------------------------------------------------------------
# DEF sheet: season-to-date defensive totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 387
$ws.Range("F2").Value = 122
$ws.Range("G2").Value = 114
$ws.Range("H2").Value = 12
$ws.Range("J2").Value = 56
$ws.Range("L2").Value = 518
$ws.Range("M2").Value = 343
$ws.Range("O2").Value = 44
$ws.Range("Q2").Value = 966
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 343
$ws.Range("E3").Value = 62
$ws.Range("F3").Value = 191
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 47
$ws.Range("I3").Value = 128
$ws.Range("J3").Value = 119
$ws.Range("N3").Value = 27

# ---------------------------------------------------------------------
# ST sheet: season-to-date special-teams totals + Week 17 samples
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 160
$ws.Range("D2").Value = 115
$ws.Range("F2").Value = 584
$ws.Range("G2").Value = 569
$ws.Range("J2").Value = 289
$ws.Range("K2").Value = 274
$ws.Range("B4").Value = $ws.Range("B4").Text + " 66 65 61 66 57"
$ws.Range("B5").Value = $ws.Range("B5").Text + " 21 18 26 20 13"
$ws.Range("B6").Value = $ws.Range("B6").Text + " 28 7"
$ws.Range("D3").Value = $ws.Range("D3").Text + " 52 28 36 43"
$ws.Range("D4").Value = $ws.Range("D4").Text + " 15 0 8 0"
$ws.Range("D5").Value = $ws.Range("D5").Text + " 0 9 0 18 11 0"

# ---------------------------------------------------------------------
# TURNS sheet: season-to-date turnover totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B2").Value = 13
$ws.Range("C2").Value = 10
$ws.Range("E2").Value = 17
$ws.Range("D3").Value = 15

# ---------------------------------------------------------------------
# PEN sheet: season-to-date penalty totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 36
$ws.Range("B3").Value = 38
$ws.Range("D4").Value = 26
